$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells (column D) are stored as text in this sheet even though many of
# them look like plain numbers (e.g. "26.80"). A leading apostrophe forces
# Excel to keep the assigned value as text instead of silently coercing it to
# a Number (which would corrupt values with trailing zeros, e.g. 26.80 -> 26.8,
# or re-interpret multi-dot figures like "62.518.87").
$ws.Range("D2").Value = '''62.518.87'
$ws.Range("E2").Value = '  +4.31%  '
$ws.Range("D3").Value = '''3.335.24'
$ws.Range("E3").Value = '  +4.27%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''559.02'
$ws.Range("E5").Value = '  +4.30%  '
$ws.Range("D6").Value = '''151.39'
$ws.Range("E6").Value = '  +4.33%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '''3.339.81'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").Value = '''0.535'
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("D10").Value = '''7.38'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("E11").Value = '  +3.65%  '
$ws.Range("D12").Value = '''0.433'
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '''3.913.78'
$ws.Range("E13").Value = '  +4.23%  '
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '''26.74'
$ws.Range("E15").Value = '  +2.88%  '
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("D17").Value = '''62.534.68'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").Value = '''3.339.88'
$ws.Range("E18").Value = '  +4.64%  '
$ws.Range("D19").Value = '''6.32'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").Value = '''13.76'
$ws.Range("E20").Value = '  +4.52%  '
$ws.Range("D21").Value = '''8.33'
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").Value = '''382.81'
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = '''0.532'
$ws.Range("E24").Value = '  +1.44%  '
$ws.Range("D25").Value = '''69.91'
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("D26").Value = '''0.177'
$ws.Range("E26").Value = '  +4.74%  '
$ws.Range("D27").Value = '''8.96'
$ws.Range("E27").Value = '  +2.19%  '
$ws.Range("D29").Value = '''0.0₃0941'
$ws.Range("E29").Value = '  +5.02%  '
$ws.Range("E30").Value = '  +5.89%  '
$ws.Range("D31").Value = '''1.97'
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''22.85'
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").Value = '''5.56'
$ws.Range("E33").Value = '  +2.49%  '
$ws.Range("E34").Value = '  +7.31%  '
$ws.Range("E35").Value = '  +0.20%  '
$ws.Range("D36").Value = '''159.52'
$ws.Range("E36").Value = '  +1.61%  '
$ws.Range("D37").Value = '''1.46'
$ws.Range("E37").Value = '  +8.60%  '
$ws.Range("D38").Value = '''1.88'
$ws.Range("E38").Value = '  +12.45%  '
$ws.Range("D39").Value = '''26.80'
$ws.Range("E39").Value = '  +5.04%  '
$ws.Range("D40").Value = '''0.0735'
$ws.Range("E40").Value = '  +4.49%  '
$ws.Range("D41").Value = '''2.795.79'
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("E42").Value = '  +6.40%  '
$ws.Range("D43").Value = '''40.41'
$ws.Range("E43").Value = '  +1.58%  '
$ws.Range("D44").Value = '''0.741'
$ws.Range("E44").Value = '  +3.44%  '
$ws.Range("D45").Value = '''4.23'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  +4.41%  '
$ws.Range("D47").Value = '''3.381.11'
$ws.Range("E47").Value = '  +4.25%  '
$ws.Range("D48").Value = '''21.79'
$ws.Range("E48").Value = '  +5.90%  '
$ws.Range("E49").Value = '  -1.48%  '
$ws.Range("D50").Value = '''6.28'
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("D51").Value = '''285.34'
$ws.Range("E51").Value = '  +5.71%  '